$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 409.8
$ws.Range("I12").Value = 399
$ws.Range("J12").Value = 412.5
$ws.Range("K12").Value = 399
$ws.Range("L12").Value = 412.5
$ws.Range("M12").Value = -229
$ws.Range("N12").Value = -752.5

$ws.Range("H33").Value = 277.5
$ws.Range("I33").Value = 208.57143
$ws.Range("J33").Value = 567
$ws.Range("K33").Value = 208.57143
$ws.Range("L33").Value = 567
$ws.Range("M33").Value = 20.42857000000001
$ws.Range("N33").Value = -1025

$ws.Range("H51").Value = 2819.1
$ws.Range("I51").Value = 2332.3333
$ws.Range("J51").Value = 3027.7144
$ws.Range("K51").Value = 2332.3333
$ws.Range("L51").Value = 3027.7144
$ws.Range("M51").Value = -1848.3333
$ws.Range("N51").Value = -3995.7144

$ws.Range("H129").Value = 831.5909
$ws.Range("J129").Value = 912.30554
$ws.Range("L129").Value = 2736.91662
$ws.Range("N129").Value = -12736.91662

$ws.Range("H138").Value = 597248.75
$ws.Range("J138").Value = 764481.8
$ws.Range("L138").Value = 2293445.4
$ws.Range("N138").Value = -2303725.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 978
$ws.Range("I2").Value = 817
$ws.Range("J2").Value = 1300
$ws.Range("K2").Value = 817
$ws.Range("L2").Value = 1300
$ws.Range("M2").Value = -704
$ws.Range("N2").Value = -1526

$ws.Range("H32").Value = 6357.759
$ws.Range("I32").Value = 6357.759
$ws.Range("K32").Value = 6357.759
$ws.Range("M32").Value = -6070.759

$ws.Range("H74").Value = 2530.9167
$ws.Range("I74").Value = 2108.875
$ws.Range("J74").Value = 3375
$ws.Range("K74").Value = 2108.875
$ws.Range("L74").Value = 3375
$ws.Range("M74").Value = -1234.875
$ws.Range("N74").Value = -5123

$ws.Range("H77").Value = 2530.9167
$ws.Range("I77").Value = 2108.875
$ws.Range("J77").Value = 3375
$ws.Range("K77").Value = 10544.375
$ws.Range("L77").Value = 16875
$ws.Range("M77").Value = -6176.375
$ws.Range("N77").Value = -25611

$ws.Range("H101").Value = 32249.75
$ws.Range("J101").Value = 32249.75
$ws.Range("L101").Value = 32249.75
$ws.Range("N101").Value = -38739.75

$ws.Range("H116").Value = 978
$ws.Range("I116").Value = 817
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 817
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = 1477
$ws.Range("N116").Value = -5888

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 978
$ws.Range("I3").Value = 817
$ws.Range("J3").Value = 1300
$ws.Range("K3").Value = 817
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = -703
$ws.Range("N3").Value = -1528

$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("N95").Value = -40492

$ws.Range("H104").Value = 65000
$ws.Range("J104").Value = 65000
$ws.Range("L104").Value = 65000
$ws.Range("N104").Value = -71988

$ws.Range("H105").Value = 125001700
$ws.Range("I105").Value = 142858800
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 142858800
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -142857053
$ws.Range("N105").Value = -5493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1280426.4
$ws.Range("J132").Value = 1280426.4
$ws.Range("L132").Value = 1280426.4
$ws.Range("N132").Value = -1290546.4

$ws.Range("H134").Value = 4814.448
$ws.Range("I134").Value = 937.94446
$ws.Range("K134").Value = 2813.83338
$ws.Range("M134").Value = -278.83338

$ws.Range("I32").Value = 1766.6666
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 1766.6666
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -1450.6666
$ws.Range("N32").Value = -10632

$ws.Range("H35").Value = 1025
$ws.Range("I35").Value = 1025
$ws.Range("K35").Value = 1025
$ws.Range("M35").Value = -731

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null

$ws.Range("H41").Value = 24566
$ws.Range("J41").Value = 24566
$ws.Range("L41").Value = 24566
$ws.Range("N41").Value = -25422

$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1407

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null

$ws.Range("H58").Value = 1689.7727
$ws.Range("I58").Value = 1390.9286
$ws.Range("K58").Value = 1390.9286
$ws.Range("M58").Value = -1187.9286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2038.4
$ws.Range("I132").Value = 1382.3077
$ws.Range("J132").Value = 3256.8572
$ws.Range("K132").Value = 4146.9231
$ws.Range("L132").Value = 9770.571599999999
$ws.Range("M132").Value = -1616.9231
$ws.Range("N132").Value = -14830.5716

$ws.Range("H134").Value = 14707463
$ws.Range("I134").Value = 1605.5172
$ws.Range("K134").Value = 4816.5516
$ws.Range("M134").Value = -2281.5516

$ws.Range("H136").Value = 1689.7727
$ws.Range("I136").Value = 1390.9286
$ws.Range("K136").Value = 4172.7858
$ws.Range("M136").Value = -1622.7858

$ws.Range("H4").Value = 3027712.5
$ws.Range("J4").Value = 2797722
$ws.Range("L4").Value = 8393166
$ws.Range("N4").Value = -8393390

$ws.Range("H5").Value = 1198
$ws.Range("I5").Value = 1354.8422
$ws.Range("J5").Value = 602
$ws.Range("K5").Value = 4064.5266
$ws.Range("L5").Value = 1806
$ws.Range("M5").Value = -3952.5266
$ws.Range("N5").Value = -2030

$ws.Range("H40").Value = 340.7143
$ws.Range("I40").Value = 117.5
$ws.Range("J40").Value = 430
$ws.Range("K40").Value = 470
$ws.Range("L40").Value = 1720
$ws.Range("M40").Value = -401
$ws.Range("N40").Value = -1858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 17860422
$ws.Range("J131").Value = 3914.8696
$ws.Range("L131").Value = 11744.6088
$ws.Range("N131").Value = -21824.6088

$ws.Range("H135").Value = 1198
$ws.Range("I135").Value = 1354.8422
$ws.Range("J135").Value = 602
$ws.Range("K135").Value = 12193.5798
$ws.Range("L135").Value = 5418
$ws.Range("M135").Value = -9658.5798
$ws.Range("N135").Value = -10488

$ws.Range("H5").Value = 15000
$ws.Range("J5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("N5").Value = -15224

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1800000
$ws.Range("J2").Value = 4400000
$ws.Range("L2").Value = 4400000
$ws.Range("N2").Value = -4400224

$ws.Range("H22").Value = 736.1667
$ws.Range("I22").Value = 477.6
$ws.Range("J22").Value = 920.8570999999999
$ws.Range("K22").Value = 477.6
$ws.Range("L22").Value = 920.8570999999999
$ws.Range("M22").Value = -182.6
$ws.Range("N22").Value = -1510.8571

$ws.Range("H27").Value = 736.1667
$ws.Range("I27").Value = 477.6
$ws.Range("J27").Value = 920.8570999999999
$ws.Range("K27").Value = 477.6
$ws.Range("L27").Value = 920.8570999999999
$ws.Range("M27").Value = -370.6
$ws.Range("N27").Value = -1134.8571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3500
$ws.Range("J2").Value = 3500
$ws.Range("L2").Value = 3500
$ws.Range("N2").Value = -3724

$ws.Range("H107").Value = 540.6667
$ws.Range("I107").Value = 441.81818
$ws.Range("K107").Value = 1325.45454
$ws.Range("M107").Value = 594.54546
